# "Add files via upload" -- refreshed figures for the Sector sheet.
#
# Net content changes applied here:
#   - Header G1 relabelled "Juntas de Freguesia" -> "Delegações CML"
#     (the shared-string table is reshaped by Excel automatically: the
#     now-unused "Juntas de Freguesia" entry drops out and the new label
#     is appended at the end).
#   - Updated figures for 2017 (row 6) and 2018 (row 7).
#   - Active selection moved from H17 to C17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the "Juntas de Freguesia" column header.
$ws.Range("G1").Value = "Delegações CML"

# 2017 row: Instituições + Total revised upward.
$ws.Range("D6").Value = 2216.8420075000004
$ws.Range("I6").Value = 57261.557134281815

# 2018 row: Instituições revised upward, Juntas de Freguesia/Delegações CML
# revised downward by the same amount (Total unchanged).
$ws.Range("D7").Value = 2097.1594574999999
$ws.Range("G7").Value = 1435.57799

# Restore the saved selection to C17.
$ws.Range("C17").Select()
